$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "bleu" status label with "noir"
$used = $ws.UsedRange
$labelCol = $used.Find("bleu")
while ($labelCol -ne $null) {
    $labelCol.Value = "noir"
    $labelCol = $used.Find("bleu")
}

# Replace the old "résultat et / ou publication posté dans les 36 mois" status name
$nameCol = $used.Find("résultat et / ou publication posté dans les 36 mois")
while ($nameCol -ne $null) {
    $nameCol.Value = "résultat postés ou publiés dans les 36 mois"
    $nameCol = $used.Find("résultat et / ou publication posté dans les 36 mois")
}

# Replace the old "pas de résultat ni de publication" status name
$noResCol = $used.Find("pas de résultat ni de publication")
while ($noResCol -ne $null) {
    $noResCol.Value = "pas de résultat postés ni publiés"
    $noResCol = $used.Find("pas de résultat ni de publication")
}
